$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Compiled")
Write-Host $ws.Name
